$d = $word.ActiveDocument

# --- Step 1: split "Version" into two runs "Versi" + "on" ---------------
# A temporary bookmark forces the engine to keep the text split into
# separate <w:r> elements instead of re-merging them into one run.
$splitPoint = $d.Range(5, 5)
$d.Bookmarks.Add("__splitmark__", $splitPoint)
$d.Bookmarks("__splitmark__").Delete()

# --- Step 2: " 1." -> " 2" (drop the trailing period here) --------------
$r = $d.Range(8, 10)
$r.Text = "2"

# --- Step 3: append a new run containing "." right after the _GoBack ---
#             bookmark (still before the paragraph mark) ----------------
$bm = $d.Bookmarks("_GoBack")
$after = $d.Range($bm.End, $bm.End)
$after.InsertAfter(".")
